# Pasajeros_AICM.xlsx – monthly update
# Adds the December-2024 (Dic.) CDMX passenger figure as a new first data
# row inside "Tabla3" (row 6 of the sheet), pushing the existing rows down
# by one, grows the table/autofilter range accordingly, and refreshes the
# "Actualización" footer note from "Noviembre 2024" to "Diciembre 2024".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a brand-new blank row above the current first data row (row 6,
#    right below the "Año/Mes/CDMX." header in row 5). Everything at/after
#    row 6 shifts down by one row, and the sheet dimension grows with it.
$ws.Rows(6).Insert()

# 2. Copy the banding/number formatting from the row that is now directly
#    below the new blank row (row 8, formerly row 7) so the new row 6 picks
#    up the same alternating-stripe style used by the rest of the table
#    instead of a generic default style.
$ws.Range("B8:D8").Copy()
$ws.Range("B6:D6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3. Fill in the new December 2024 figures.
$ws.Range("B6").Value = 2024
$ws.Range("C6").Value = "Dic."
$ws.Range("D6").Value = 3901.751

# 4. Grow the Excel Table ("Tabla3") / its AutoFilter range by one row so it
#    keeps covering the data (B5:D88 -> B5:D89).
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B5:D89"))

# 5. Update the "Actualización: ..." footer note to reflect the new month.
$ws.Range("B90").Value = "Actualización: Diciembre 2024."
